# Rename header labels on the existing sheets
$wb = $excel.ActiveWorkbook

$weekly = $wb.Worksheets.Item("Weekly Quantity")
$weekly.Range("B1").Value = "Weekly_PO_Qty"

$monthly = $wb.Worksheets.Item("Monthly Trend")
$monthly.Range("B1").Value = "Monthly_PO_Qty"

# Add the new "PO Forecast" sheet after the last existing sheet
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add($null, $lastSheet)
$newSheet.Name = "PO Forecast"

# Header row
$newSheet.Range("A1").Value = "ds"
$newSheet.Range("B1").Value = "PO_Forecast"
$newSheet.Range("C1").Value = "yhat_lower"
$newSheet.Range("D1").Value = "yhat_upper"

# Match header formatting (bold, centered, bordered) used on the other sheets
$weekly.Range("A1:B1").Copy()
$newSheet.Range("A1:B1").PasteSpecial(-4122)
$newSheet.Range("C1:D1").PasteSpecial(-4122)

# Data rows
$newSheet.Range("A2").Value = 45417.99999999999
$newSheet.Range("B2").Value = 52
$newSheet.Range("C2").Value = -35.83131171093311
$newSheet.Range("D2").Value = 139.3798570226398
$newSheet.Range("A3").Value = 45445.99999999999
$newSheet.Range("B3").Value = 64
$newSheet.Range("C3").Value = -19.65285728261248
$newSheet.Range("D3").Value = 150.6998602795921
$newSheet.Range("A4").Value = 45459.99999999999
$newSheet.Range("B4").Value = 70
$newSheet.Range("C4").Value = -18.52941909121531
$newSheet.Range("D4").Value = 160.5739280100676
$newSheet.Range("A5").Value = 45466.99999999999
$newSheet.Range("B5").Value = 73
$newSheet.Range("C5").Value = -16.24407806149293
$newSheet.Range("D5").Value = 164.432805005101
$newSheet.Range("A6").Value = 45480.99999999999
$newSheet.Range("B6").Value = 80
$newSheet.Range("C6").Value = -6.517056065891336
$newSheet.Range("D6").Value = 168.9400643267328
$newSheet.Range("A7").Value = 45487.99999999999
$newSheet.Range("B7").Value = 83
$newSheet.Range("C7").Value = -1.739768495444172
$newSheet.Range("D7").Value = 176.2008489798099
$newSheet.Range("A8").Value = 45501.99999999999
$newSheet.Range("B8").Value = 89
$newSheet.Range("C8").Value = 1.086087795739941
$newSheet.Range("D8").Value = 178.1808061950969
$newSheet.Range("A9").Value = 45508.99999999999
$newSheet.Range("B9").Value = 92
$newSheet.Range("C9").Value = 4.119748871822011
$newSheet.Range("D9").Value = 182.9529088406698
$newSheet.Range("A10").Value = 45515.99999999999
$newSheet.Range("B10").Value = 95
$newSheet.Range("C10").Value = 13.66341352949337
$newSheet.Range("D10").Value = 187.9842242242218
$newSheet.Range("A11").Value = 45522.99999999999
$newSheet.Range("B11").Value = 98
$newSheet.Range("C11").Value = 8.815367623255252
$newSheet.Range("D11").Value = 183.6917172900646
$newSheet.Range("A12").Value = 45529.99999999999
$newSheet.Range("B12").Value = 102
$newSheet.Range("C12").Value = 9.520096542609277
$newSheet.Range("D12").Value = 195.0360203426382
$newSheet.Range("A13").Value = 45557.99999999999
$newSheet.Range("B13").Value = 114
$newSheet.Range("C13").Value = 23.5324784825537
$newSheet.Range("D13").Value = 203.4493808229534
$newSheet.Range("A14").Value = 45571.99999999999
$newSheet.Range("B14").Value = 120
$newSheet.Range("C14").Value = 30.30156017078881
$newSheet.Range("D14").Value = 217.5524498513485
$newSheet.Range("A15").Value = 45585.99999999999
$newSheet.Range("B15").Value = 127
$newSheet.Range("C15").Value = 37.89675790558262
$newSheet.Range("D15").Value = 214.3035699845566
$newSheet.Range("A16").Value = 45592.99999999999
$newSheet.Range("B16").Value = 130
$newSheet.Range("C16").Value = 41.3864769767954
$newSheet.Range("D16").Value = 216.9596067381317
$newSheet.Range("A17").Value = 45599.99999999999
$newSheet.Range("B17").Value = 133
$newSheet.Range("C17").Value = 39.41876717334407
$newSheet.Range("D17").Value = 221.9236004762817
$newSheet.Range("A18").Value = 45606.99999999999
$newSheet.Range("B18").Value = 136
$newSheet.Range("C18").Value = 49.76067504495989
$newSheet.Range("D18").Value = 223.7682436065298
$newSheet.Range("A19").Value = 45613.99999999999
$newSheet.Range("B19").Value = 139
$newSheet.Range("C19").Value = 52.31542985957697
$newSheet.Range("D19").Value = 231.6810183958946
$newSheet.Range("A20").Value = 45620.99999999999
$newSheet.Range("B20").Value = 142
$newSheet.Range("C20").Value = 54.68512507938642
$newSheet.Range("D20").Value = 238.4597635220683
$newSheet.Range("A21").Value = 45627.99999999999
$newSheet.Range("B21").Value = 145
$newSheet.Range("C21").Value = 56.69159723966413
$newSheet.Range("D21").Value = 238.5587761159526
$newSheet.Range("A22").Value = 45634.99999999999
$newSheet.Range("B22").Value = 148
$newSheet.Range("C22").Value = 60.59624491748605
$newSheet.Range("D22").Value = 238.8396342502056
$newSheet.Range("A23").Value = 45641.99999999999
$newSheet.Range("B23").Value = 152
$newSheet.Range("C23").Value = 62.48644901675252
$newSheet.Range("D23").Value = 240.4903417920573
$newSheet.Range("A24").Value = 45648.99999999999
$newSheet.Range("B24").Value = 155
$newSheet.Range("C24").Value = 71.21994636610181
$newSheet.Range("D24").Value = 246.2343365628868
$newSheet.Range("A25").Value = 45655.99999999999
$newSheet.Range("B25").Value = 158
$newSheet.Range("C25").Value = 68.57914777299979
$newSheet.Range("D25").Value = 243.1387676182925

# Match the date-formatted style used for column A on the other sheets
$weekly.Range("A2").Copy()
$newSheet.Range("A2:A25").PasteSpecial(-4122)
